# Append new translation rows to the "Import" sheet (cascade create/validation
# labels for vendor/atomizer/build), mirroring rows already present there.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# New rows to append, in order: (key, translation)
$rows = @(
    @("lab.vendor.name.label.tooltip", "Při zadávání jména výrobce se prosím snažte držet přesného názvu, včetně velikosti písmen a cizích znaků."),
    @("lab.vendor.name.label", "Jméno výrobce"),
    @("lab.vendor.create.submit", "Vytvořit"),
    @("lab.atomizer.create.submit", "Vytvořit"),
    @("lab.vendor.name.label.required", "Jméno výrobce je povinné"),
    @("lab.atomizer.name.label.required", "Jméno atomizéru je povinné"),
    @("lab.build.name.label.required", "Jméno buildu je povinné"),
    @("lab.build.create.submit", "Vytvořit"),
    @("lab.vendor.create.success", "Výrobce [{{data.name}}] byl úspěšně vytvořen."),
    @("lab.atomizer.create.success", "Atomizér [{{data.name}}] byl úspěšně vytvořen.")
)

$lastRow = 218
$startRow = $lastRow + 1
$endRow = $lastRow + $rows.Length

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $pair = $rows[$i]

    $null = $ws.Range("A" + $lastRow + ":C" + $lastRow).Copy()
    $null = $ws.Range("A" + $r + ":C" + $r).PasteSpecial(-4122)

    $ws.Range("A" + $r).Value = "cs"
    $ws.Range("B" + $r).Value = $pair[0]
    $ws.Range("C" + $r).Value = $pair[1]
}

$null = $ws.Range("B" + $startRow).Select()
